# Update "想去人数" (want-to-go count) values in column F for rows 2-14
# on both the "展览" and "全部类型" worksheets.

$wb = $excel.ActiveWorkbook

$newValues = @{
    2  = 1760
    3  = 243
    4  = 226
    5  = 7341
    6  = 499
    7  = 508
    8  = 70
    9  = 19
    10 = 9029
    11 = 2381
    12 = 293
    13 = 9161
    14 = 10457
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $newValues.Keys) {
        $ws.Range("F$row").Value = $newValues[$row]
    }
}
